$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 124
$ws.Range("H124").Value = 45962.8
$ws.Range("J124").Value = 45962.8
$ws.Range("L124").Value = 45962.8
$ws.Range("N124").Value = -55782.8

# Row 128
$ws.Range("H128").Value = 52177
$ws.Range("J128").Value = 52177
$ws.Range("L128").Value = 52177
$ws.Range("N128").Value = -62137

# Row 129
$ws.Range("H129").Value = 318140.44
$ws.Range("J129").Value = 2500
$ws.Range("L129").Value = 7500
$ws.Range("N129").Value = -17500

# Row 130
$ws.Range("H130").Value = 54996
$ws.Range("J130").Value = 54996
$ws.Range("L130").Value = 54996
$ws.Range("N130").Value = -65036

# Row 138
$ws.Range("H138").Value = 3108.5422
$ws.Range("I138").Value = 2922
$ws.Range("J138").Value = 3163.9219
$ws.Range("K138").Value = 8766
$ws.Range("L138").Value = 9491.7657
$ws.Range("M138").Value = -3626
$ws.Range("N138").Value = -19771.7657

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2074.889
$ws.Range("I2").Value = 2193.5
$ws.Range("J2").Value = 1837.6666
$ws.Range("K2").Value = 2193.5
$ws.Range("L2").Value = 1837.6666
$ws.Range("M2").Value = -2080.5
$ws.Range("N2").Value = -2063.6666

# Row 32
$ws.Range("H32").Value = 23769.297
$ws.Range("I32").Value = 23769.297
$ws.Range("K32").Value = 23769.297
$ws.Range("M32").Value = -23482.297

# Row 74
$ws.Range("H74").Value = 2308.1724
$ws.Range("I74").Value = 1954.5625
$ws.Range("J74").Value = 2743.3845
$ws.Range("K74").Value = 1954.5625
$ws.Range("L74").Value = 2743.3845
$ws.Range("M74").Value = -1080.5625
$ws.Range("N74").Value = -4491.3845

# Row 77
$ws.Range("H77").Value = 2308.1724
$ws.Range("I77").Value = 1954.5625
$ws.Range("J77").Value = 2743.3845
$ws.Range("K77").Value = 9772.8125
$ws.Range("L77").Value = 13716.9225
$ws.Range("M77").Value = -5404.8125
$ws.Range("N77").Value = -22452.9225

# Row 114
$ws.Range("H114").Value = 34111.5
$ws.Range("J114").Value = 34111.5
$ws.Range("L114").Value = 34111.5
$ws.Range("N114").Value = -42789.5

# Row 116
$ws.Range("H116").Value = 2074.889
$ws.Range("I116").Value = 2193.5
$ws.Range("J116").Value = 1837.6666
$ws.Range("K116").Value = 2193.5
$ws.Range("L116").Value = 1837.6666
$ws.Range("M116").Value = 100.5
$ws.Range("N116").Value = -6425.6666

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Row 130
$ws.Range("H130").Value = 46916.668
$ws.Range("J130").Value = 46916.668
$ws.Range("L130").Value = 46916.668
$ws.Range("N130").Value = -56956.668

# Row 131
$ws.Range("H131").Value = 50357
$ws.Range("J131").Value = 50357
$ws.Range("L131").Value = 50357
$ws.Range("N131").Value = -60437

# Row 133
$ws.Range("H133").Value = 32113.143
$ws.Range("J133").Value = 32113.143
$ws.Range("L133").Value = 32113.143
$ws.Range("N133").Value = -37173.143

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2074.889
$ws.Range("I3").Value = 2193.5
$ws.Range("J3").Value = 1837.6666
$ws.Range("K3").Value = 2193.5
$ws.Range("L3").Value = 1837.6666
$ws.Range("M3").Value = -2079.5
$ws.Range("N3").Value = -2065.6666

# Row 92
$ws.Range("H92").Value = 44071
$ws.Range("J92").Value = 44071
$ws.Range("L92").Value = 44071
$ws.Range("N92").Value = -49063

# Row 100
$ws.Range("H100").Value = 45643
$ws.Range("J100").Value = 45643
$ws.Range("L100").Value = 45643
$ws.Range("N100").Value = -47807

# Row 124
$ws.Range("H124").Value = 49850.668
$ws.Range("J124").Value = 49850.668
$ws.Range("L124").Value = 49850.668
$ws.Range("N124").Value = -59670.668

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# Row 130
$ws.Range("H130").Value = 45309.75
$ws.Range("J130").Value = 45309.75
$ws.Range("L130").Value = 45309.75
$ws.Range("N130").Value = -55349.75

# Row 135
$ws.Range("H135").Value = 44544.125
$ws.Range("J135").Value = 44544.125
$ws.Range("L135").Value = 44544.125
$ws.Range("N135").Value = -54684.125

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1639.6154
$ws.Range("I16").Value = 1412.75
$ws.Range("J16").Value = 2002.6
$ws.Range("K16").Value = 1412.75
$ws.Range("L16").Value = 2002.6
$ws.Range("M16").Value = -1125.75
$ws.Range("N16").Value = -2576.6

# Row 20
$ws.Range("H20").Value = 49772
$ws.Range("J20").Value = 49772
$ws.Range("L20").Value = 49772
$ws.Range("N20").Value = -50244

# Row 30
$ws.Range("H30").Value = 49772
$ws.Range("J30").Value = 49772
$ws.Range("L30").Value = 49772
$ws.Range("N30").Value = -49954

# Row 87
$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27372

# Row 90
$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -86856

# Row 95
$ws.Range("H95").Value = 200000
$ws.Range("J95").Value = 200000
$ws.Range("L95").Value = 200000
$ws.Range("N95").Value = -205492

# Row 99
$ws.Range("H99").Value = 2214.0833
$ws.Range("I99").Value = 2308.8572
$ws.Range("J99").Value = 2081.4
$ws.Range("K99").Value = 2308.8572
$ws.Range("L99").Value = 2081.4
$ws.Range("M99").Value = -810.8571999999999
$ws.Range("N99").Value = -5077.4

# Row 113
$ws.Range("H113").Value = 1639.6154
$ws.Range("I113").Value = 1412.75
$ws.Range("J113").Value = 2002.6
$ws.Range("K113").Value = 1412.75
$ws.Range("L113").Value = 2002.6
$ws.Range("M113").Value = 757.25
$ws.Range("N113").Value = -6342.6

# Row 126
$ws.Range("H126").Value = 2214.0833
$ws.Range("I126").Value = 2308.8572
$ws.Range("J126").Value = 2081.4
$ws.Range("K126").Value = 6926.571599999999
$ws.Range("L126").Value = 6244.200000000001
$ws.Range("M126").Value = -4456.571599999999
$ws.Range("N126").Value = -11184.2

# Row 128
$ws.Range("H128").Value = 49772
$ws.Range("J128").Value = 49772
$ws.Range("L128").Value = 49772
$ws.Range("N128").Value = -59732

# Row 134
$ws.Range("H134").Value = 2367.4348
$ws.Range("I134").Value = 1483.625
$ws.Range("J134").Value = 4387.5713
$ws.Range("K134").Value = 4450.875
$ws.Range("L134").Value = 13162.7139
$ws.Range("M134").Value = -1915.875
$ws.Range("N134").Value = -18232.7139

# Row 135
$ws.Range("H135").Value = 55704.844
$ws.Range("J135").Value = 55704.844
$ws.Range("L135").Value = 55704.844
$ws.Range("N135").Value = -65844.844

$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 6836.4614
$ws.Range("I56").Value = 6836.4614
$ws.Range("K56").Value = 6836.4614
$ws.Range("M56").Value = -6306.4614

# Row 131
$ws.Range("H131").Value = 31212.457
$ws.Range("I131").Value = 9634.454
$ws.Range("J131").Value = 35235.473
$ws.Range("K131").Value = 28903.362
$ws.Range("L131").Value = 105706.419
$ws.Range("M131").Value = -23863.362
$ws.Range("N131").Value = -115786.419

$ws = $wb.Worksheets.Item("GSM")
# Row 87
$ws.Range("H87").Value = 25625
$ws.Range("J87").Value = 25625
$ws.Range("L87").Value = 25625
$ws.Range("N87").Value = -28121

# Row 90
$ws.Range("H90").Value = 25625
$ws.Range("J90").Value = 25625
$ws.Range("L90").Value = 76875
$ws.Range("N90").Value = -89355

# Row 122
$ws.Range("H122").Value = 1214
$ws.Range("I122").Value = 1155.5454
$ws.Range("J122").Value = 1374.75
$ws.Range("K122").Value = 3466.6362
$ws.Range("L122").Value = 4124.25
$ws.Range("M122").Value = -1016.6362
$ws.Range("N122").Value = -9024.25

# Row 130
$ws.Range("H130").Value = 49678
$ws.Range("J130").Value = 49678
$ws.Range("L130").Value = 49678
$ws.Range("N130").Value = -59718

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 16669748
$ws.Range("I82").Value = 4501
$ws.Range("J82").Value = 27779914
$ws.Range("K82").Value = 4501
$ws.Range("L82").Value = 27779914
$ws.Range("M82").Value = -4140
$ws.Range("N82").Value = -27780636

# Row 85
$ws.Range("H85").Value = 16669748
$ws.Range("I85").Value = 4501
$ws.Range("J85").Value = 27779914
$ws.Range("K85").Value = 4501
$ws.Range("L85").Value = 27779914
$ws.Range("M85").Value = -3253
$ws.Range("N85").Value = -27782410

# Row 130
$ws.Range("H130").Value = 41227.668
$ws.Range("J130").Value = 41227.668
$ws.Range("L130").Value = 41227.668
$ws.Range("N130").Value = -51267.668

# Row 135
$ws.Range("H135").Value = 78000
$ws.Range("J135").Value = 78000
$ws.Range("L135").Value = 78000
$ws.Range("N135").Value = -88140

$ws = $wb.Worksheets.Item("WVR")
# Row 131
$ws.Range("H131").Value = 46521.2
$ws.Range("J131").Value = 46521.2
$ws.Range("L131").Value = 46521.2
$ws.Range("N131").Value = -56601.2

# Row 135
$ws.Range("H135").Value = 53057
$ws.Range("J135").Value = 53057
$ws.Range("L135").Value = 53057
$ws.Range("N135").Value = -63197
